# Generate Report for Handoff
# Adds a new "d19ebcd2-536c-4bd6-8fd7-c3116e265127" handoff entry ahead of
# the existing "ec4d2e15-c547-4e85-9055-752a790798d2" entry on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Push the existing data row (currently row 2) down to row 3, carrying
# its formatting (hyperlink style) along with it.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new handoff's summary info.
$ws.Range("A2").Value2 = "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md"
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "2016-28-17 18:28:28"
$ws.Range("A2").Style = "HyperLink"

# Rebuild hyperlinks for both rows (row insert does not relocate them).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/d19ebcd2-536c-4bd6-8fd7-c3116e265127.md", [Type]::Missing, [Type]::Missing, "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/ec4d2e15-c547-4e85-9055-752a790798d2.md", [Type]::Missing, [Type]::Missing, "ec4d2e15-c547-4e85-9055-752a790798d2.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value2 = "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.zh-cn.xlf"
$ws.Range("E2").Value2 = "2016-03-17 18:28:25"
$ws.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws.Range("I2").Value2 = "Include"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Style = "HyperLink"
$ws.Range("D2").Style = "HyperLink"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/d19ebcd2-536c-4bd6-8fd7-c3116e265127.md", [Type]::Missing, [Type]::Missing, "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/d19ebcd2-536c-4bd6-8fd7-c3116e265127.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be47c9d223d37ff722eaa6a878fa1b29ab43bd01/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/ec4d2e15-c547-4e85-9055-752a790798d2.md", [Type]::Missing, [Type]::Missing, "ec4d2e15-c547-4e85-9055-752a790798d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/ec4d2e15-c547-4e85-9055-752a790798d2.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be47c9d223d37ff722eaa6a878fa1b29ab43bd01/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ec4d2e15-c547-4e85-9055-752a790798d2.2d47fc933ea65516dfaa874554d219704769d7a1.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ec4d2e15-c547-4e85-9055-752a790798d2.2d47fc933ea65516dfaa874554d219704769d7a1.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value2 = "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.de-de.xlf"
$ws.Range("E2").Value2 = "2016-03-17 18:28:28"
$ws.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws.Range("I2").Value2 = "Include"

$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Style = "HyperLink"
$ws.Range("D2").Style = "HyperLink"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/d19ebcd2-536c-4bd6-8fd7-c3116e265127.md", [Type]::Missing, [Type]::Missing, "d19ebcd2-536c-4bd6-8fd7-c3116e265127.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/d19ebcd2-536c-4bd6-8fd7-c3116e265127.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a987a5c78693b936718abc636f42230442f5306/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.de-de.xlf", [Type]::Missing, [Type]::Missing, "d19ebcd2-536c-4bd6-8fd7-c3116e265127.372dffa7aeb20c61161314d5eddf5f11369d8bbf.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/ec4d2e15-c547-4e85-9055-752a790798d2.md", [Type]::Missing, [Type]::Missing, "ec4d2e15-c547-4e85-9055-752a790798d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/32079157c2bd503f3ba918f76435ee6ca534871b/e2e/ec4d2e15-c547-4e85-9055-752a790798d2.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a987a5c78693b936718abc636f42230442f5306/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ec4d2e15-c547-4e85-9055-752a790798d2.2d47fc933ea65516dfaa874554d219704769d7a1.de-de.xlf", [Type]::Missing, [Type]::Missing, "ec4d2e15-c547-4e85-9055-752a790798d2.2d47fc933ea65516dfaa874554d219704769d7a1.de-de.xlf") | Out-Null
